$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's used range
$lastRow = $ws.UsedRange.Rows.Count

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2026-02-28 -> 2026-03-01, i.e. 46081 -> 46082) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 1
    }
}
